# Fill in the "Absent" column (H) of the consolidated attendance report.
# Absent = 1 when there was no "Real" attendance recorded (column E = 0),
# and 0 when the student was really present (column E = 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 0
